$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin Price (D) and Volume(1h) (E) columns, and for the two rows
# where ranking order changed, also update Coin (B) and Link (C).
# Price values that are numeric-looking (e.g. "1.003") are written with a
# leading apostrophe so Excel keeps them as literal text, matching the
# original inline-string cell content instead of converting them to numbers.

# Row 2
$ws.Range('D2').Value = '29.174.63'
$ws.Range('E2').Value = '  +1.84%  '

# Row 3
$ws.Range('D3').Value = '1.909.73'
$ws.Range('E3').Value = '  +2.06%  '

# Row 4
$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  -0.27%  '

# Row 5
$ws.Range('D5').Value = '''327.52'
$ws.Range('E5').Value = '  +0.50%  '

# Row 6
$ws.Range('D6').Value = '''1.003'
$ws.Range('E6').Value = '  -0.05%  '

# Row 7
$ws.Range('D7').Value = '''0.4660'
$ws.Range('E7').Value = '  +0.32%  '

# Row 8
$ws.Range('E8').Value = '  +1.08%  '

# Row 9
$ws.Range('D9').Value = '''46.96'
$ws.Range('E9').Value = '  +0.66%  '

# Row 10
$ws.Range('D10').Value = '''0.07978'
$ws.Range('E10').Value = '  +1.31%  '

# Row 11
$ws.Range('D11').Value = '''1.005'
$ws.Range('E11').Value = '  +3.07%  '

# Row 12
$ws.Range('D12').Value = '''22.30'
$ws.Range('E12').Value = '  +1.47%  '

# Row 13
$ws.Range('D13').Value = '1.926.63'
$ws.Range('E13').Value = '  -1.68%  '

# Row 14
$ws.Range('D14').Value = '''7.140'
$ws.Range('E14').Value = '  +2.07%  '

# Row 15
$ws.Range('D15').Value = '''5.796'
$ws.Range('E15').Value = '  +1.62%  '

# Row 16
$ws.Range('D16').Value = '''0.06984'

# Row 17
$ws.Range('D17').Value = '''88.74'
$ws.Range('E17').Value = '  +0.82%  '

# Row 18
$ws.Range('D18').Value = '''1.004'
$ws.Range('E18').Value = '  -0.16%  '

# Row 19
$ws.Range('E19').Value = '  +0.62%  '

# Row 20
$ws.Range('D20').Value = '''17.26'
$ws.Range('E20').Value = '  +2.53%  '

# Row 21
$ws.Range('D21').Value = '''1.004'
$ws.Range('E21').Value = '  -0.05%  '

# Row 22
$ws.Range('D22').Value = '29.182.23'
$ws.Range('E22').Value = '  +1.79%  '

# Row 23
$ws.Range('D23').Value = '''5.378'
$ws.Range('E23').Value = '  +1.72%  '

# Row 24
$ws.Range('E24').Value = '  +0.73%  '

# Row 25
$ws.Range('D25').Value = '2.154.44'
$ws.Range('E25').Value = '  -1.46%  '

# Row 26
$ws.Range('E26').Value = '  -2.88%  '

# Row 27
$ws.Range('E27').Value = '  +2.33%  '

# Row 28
$ws.Range('D28').Value = '''19.56'
$ws.Range('E28').Value = '  +1.77%  '

# Row 29
$ws.Range('D29').Value = '''5.853'
$ws.Range('E29').Value = '  +1.09%  '

# Row 30
$ws.Range('D30').Value = '''2.007'
$ws.Range('E30').Value = '  +0.98%  '

# Row 31
$ws.Range('D31').Value = '''119.66'
$ws.Range('E31').Value = '  +0.24%  '

# Row 32
$ws.Range('D32').Value = '''0.09411'
$ws.Range('E32').Value = '  +0.52%  '

# Row 33
$ws.Range('D33').Value = '''0.9250'
$ws.Range('E33').Value = '  +0.45%  '

# Row 34
$ws.Range('D34').Value = '''5.369'
$ws.Range('E34').Value = '  +1.95%  '

# Row 35
$ws.Range('D35').Value = '''1.346'
$ws.Range('E35').Value = '  +0.64%  '

# Row 36
$ws.Range('D36').Value = '''3.269'
$ws.Range('E36').Value = '  -1.82%  '

# Row 37
$ws.Range('D37').Value = '''0.05849'
$ws.Range('E37').Value = '  +0.93%  '

# Row 38
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').Value = '''8.034'
$ws.Range('E38').Value = '  +3.41%  '

# Row 39
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '''1.159'
$ws.Range('E39').Value = '  +1.04%  '

# Row 40
$ws.Range('D40').Value = '''0.02099'

# Row 41
$ws.Range('D41').Value = '''0.5760'
$ws.Range('E41').Value = '  +2.40%  '

# Row 42
$ws.Range('D42').Value = '''0.1811'
$ws.Range('E42').Value = '  +1.42%  '

# Row 43
$ws.Range('D43').Value = '''10.02'
$ws.Range('E43').Value = '  +2.56%  '

# Row 44
$ws.Range('D44').Value = '''12.05'
$ws.Range('E44').Value = '  +2.44%  '

# Row 45
$ws.Range('D45').Value = '''0.5433'
$ws.Range('E45').Value = '  +2.25%  '

# Row 46
$ws.Range('D46').Value = '''2.228'
$ws.Range('E46').Value = '  +5.75%  '

# Row 47
$ws.Range('D47').Value = '''0.07106'
$ws.Range('E47').Value = '  -1.41%  '

# Row 48
$ws.Range('D48').Value = '''1.886'
$ws.Range('E48').Value = '  +3.32%  '

# Row 49
$ws.Range('D49').Value = '''2.599'
$ws.Range('E49').Value = '  +7.27%  '

# Row 50
$ws.Range('D50').Value = '''112.44'
$ws.Range('E50').Value = '  -0.68%  '

# Row 51
$ws.Range('D51').Value = '''1.085'
$ws.Range('E51').Value = '  -5.53%  '
